$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.336.03"
$ws.Range("E2").Value = "  +3.50%  "

$ws.Range("D3").Value = "3.111.92"
$ws.Range("E3").Value = "  +1.49%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'219.74"
$ws.Range("E5").Value = "  +4.50%  "

$ws.Range("D6").Value = "'622.33"
$ws.Range("E6").Value = "  +0.63%  "

$ws.Range("D7").Value = "'0.380"
$ws.Range("E7").Value = "  +2.91%  "

$ws.Range("D8").Value = "'0.972"
$ws.Range("E8").Value = "  +21.55%  "

$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").Value = "3.107.76"
$ws.Range("E10").Value = "  +1.44%  "

$ws.Range("D11").Value = "'0.718"
$ws.Range("E11").Value = "  +21.40%  "

$ws.Range("E12").Value = "  +5.12%  "

$ws.Range("E13").Value = "  +7.12%  "

$ws.Range("D14").Value = "'34.48"
$ws.Range("E14").Value = "  +8.27%  "

$ws.Range("D15").Value = "91.175.81"
$ws.Range("E15").Value = "  +3.64%  "

$ws.Range("E16").Value = "  +2.01%  "

$ws.Range("D17").Value = "3.691.41"
$ws.Range("E17").Value = "  +1.60%  "

$ws.Range("D18").Value = "3.122.86"
$ws.Range("E18").Value = "  +1.93%  "

$ws.Range("D19").Value = "'3.72"
$ws.Range("E19").Value = "  +13.60%  "

$ws.Range("E20").Value = "  +9.74%  "

$ws.Range("D21").Value = "'14.09"
$ws.Range("E21").Value = "  +6.34%  "

$ws.Range("D22").Value = "'436.80"
$ws.Range("E22").Value = "  +3.92%  "

$ws.Range("E23").Value = "  +7.99%  "

$ws.Range("D24").Value = "'5.17"
$ws.Range("E24").Value = "  +6.28%  "

$ws.Range("D25").Value = "'6.16"
$ws.Range("E25").Value = "  +13.14%  "

$ws.Range("D26").Value = "'87.41"
$ws.Range("E26").Value = "  +6.49%  "

$ws.Range("D27").Value = "'12.18"
$ws.Range("E27").Value = "  +3.69%  "

$ws.Range("D28").Value = "3.286.24"
$ws.Range("E28").Value = "  +1.68%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("D30").Value = "'0.167"
$ws.Range("E30").Value = "  -1.36%  "

$ws.Range("D31").Value = "'9.10"
$ws.Range("E31").Value = "  +13.55%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'526.59"
$ws.Range("E32").Value = "  +3.67%  "

$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").Value = "'0.889"
$ws.Range("E33").Value = "  -18.23%  "

$ws.Range("D34").Value = "'3.73"
$ws.Range("E34").Value = "  +4.11%  "

$ws.Range("D35").Value = "'7.09"
$ws.Range("E35").Value = "  +5.35%  "

$ws.Range("E36").Value = "  +9.64%  "

$ws.Range("D37").Value = "'23.73"
$ws.Range("E37").Value = "  +6.73%  "

$ws.Range("D38").Value = "'1.85"
$ws.Range("E38").Value = "  +3.66%  "

$ws.Range("E39").Value = "  +3.27%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.0875"
$ws.Range("E40").Value = "  +27.12%  "

$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D41").Value = "'22.30"
$ws.Range("E41").Value = "  +0.37%  "

$ws.Range("E42").Value = "  -0.08%  "

$ws.Range("E43").Value = "  +15.06%  "

$ws.Range("D45").Value = "'0.380"
$ws.Range("E45").Value = "  +5.79%  "

$ws.Range("D46").Value = "'1.92"
$ws.Range("E46").Value = "  +6.64%  "

$ws.Range("D47").Value = "'146.80"
$ws.Range("E47").Value = "  -1.20%  "

$ws.Range("D48").Value = "'44.09"
$ws.Range("E48").Value = "  +1.69%  "

$ws.Range("E49").Value = "  +9.39%  "

$ws.Range("D50").Value = "'166.97"
$ws.Range("E50").Value = "  +6.16%  "

$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D51").Value = "'4.18"
$ws.Range("E51").Value = "  +6.83%  "
